$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.506.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.923.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.76%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4711"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.92%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2880"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.66%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06741"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.75%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'106.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.74%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'18.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.74%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07758"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.46%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.905.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.79%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.56%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.28%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'292.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.17%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.501.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.93%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007580"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.52%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.150.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.89%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.267"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.199"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.17%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.354"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.26%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'168.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.37%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'21.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.69%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.097"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.97%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.1069"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.41%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.42%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.168"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.62%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.35%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.28%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7430"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.07%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.97%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.02119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +7.18%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.43%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.682"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.98%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.076"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.79%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'110.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.66%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8760"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.863"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.60%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.4270"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'67.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.55%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'49.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +15.66%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'7.188"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.72%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.256"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'35.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.49%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.34%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.2468"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +10.27%  "
$ws.Range("E51").Style = "Normal"
